$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds serial 45180 (2023-09-11) for rows 2-17.
# Bump it by one day to serial 45181 (2023-09-12).
for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45180) {
        $cell.Value = 45181
    }
}
